# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style) from the last existing header cell
# (AC1) onto the three new header cells so they match the bold/bordered
# header look, then set their labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill every data row (2-65) with the team's season record.
$ws.Range("AD2:AD65").Value = 73
$ws.Range("AE2:AE65").Value = 89
$ws.Range("AF2:AF65").Value = 0
